$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.595741748809814
$ws.Range("B1").Value = 3.85477614402771
$ws.Range("C1").Value = 3.381373643875122
$ws.Range("D1").Value = 4.175801753997803
$ws.Range("E1").Value = 4.927867412567139
